$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value2 = 618.625
$ws.Range("I12").Value2 = 531.8
$ws.Range("J12").Value2 = 763.3333
$ws.Range("K12").Value2 = 531.8
$ws.Range("L12").Value2 = 763.3333
$ws.Range("M12").Value2 = -361.8
$ws.Range("N12").Value2 = -1103.3333
$ws.Range("H34").Value2 = 1009.75
$ws.Range("I34").Value2 = 1009.75
$ws.Range("K34").Value2 = 1009.75
$ws.Range("M34").Value2 = -806.75
$ws.Range("H36").Value2 = 1009.75
$ws.Range("I36").Value2 = 1009.75
$ws.Range("K36").Value2 = 1009.75
$ws.Range("M36").Value2 = -294.75
$ws.Range("H107").Value2 = 867.3143
$ws.Range("I107").Value2 = 774.2069
$ws.Range("J107").Value2 = 1317.3334
$ws.Range("K107").Value2 = 774.2069
$ws.Range("L107").Value2 = 1317.3334
$ws.Range("M107").Value2 = 1145.7931
$ws.Range("N107").Value2 = -5157.3334
$ws.Range("H108").Value2 = 63333
$ws.Range("J108").Value2 = 63333
$ws.Range("L108").Value2 = 63333
$ws.Range("N108").Value2 = -71013
$ws.Range("H110").Value2 = 32235.5
$ws.Range("J110").Value2 = 32235.5
$ws.Range("L110").Value2 = 32235.5
$ws.Range("N110").Value2 = -40415.5
$ws.Range("H113").Value2 = 6920.5
$ws.Range("I113").Value2 = 5484.1665
$ws.Range("J113").Value2 = 9075
$ws.Range("K113").Value2 = 5484.1665
$ws.Range("L113").Value2 = 9075
$ws.Range("M113").Value2 = -2230.1665
$ws.Range("N113").Value2 = -15583
$ws.Range("H132").Value2 = 6883.5293
$ws.Range("I132").Value2 = 6860.143
$ws.Range("J132").Value2 = 6992.6665
$ws.Range("K132").Value2 = 20580.429
$ws.Range("L132").Value2 = 20977.9995
$ws.Range("M132").Value2 = -18050.429
$ws.Range("N132").Value2 = -26037.9995

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 2082.2222
$ws.Range("I2").Value2 = 1330.5834
$ws.Range("K2").Value2 = 1330.5834
$ws.Range("M2").Value2 = -1217.5834
$ws.Range("H37").Value2 = 24375
$ws.Range("J37").Value2 = 24375
$ws.Range("L37").Value2 = 24375
$ws.Range("N37").Value2 = -24921
$ws.Range("H44").Value2 = 35000
$ws.Range("J44").Value2 = 35000
$ws.Range("L44").Value2 = 35000
$ws.Range("N44").Value2 = -35976
$ws.Range("H55").Value2 = 25000
$ws.Range("J55").Value2 = 25000
$ws.Range("L55").Value2 = 25000
$ws.Range("N55").Value2 = -25630
$ws.Range("H80").Value2 = 40000
$ws.Range("J80").Value2 = 40000
$ws.Range("L80").Value2 = 40000
$ws.Range("N80").Value2 = -41996
$ws.Range("H83").Value2 = 40000
$ws.Range("J83").Value2 = 40000
$ws.Range("L83").Value2 = 120000
$ws.Range("N83").Value2 = -129984
$ws.Range("H110").Value2 = 3566.5
$ws.Range("I110").Value2 = 2010.6666
$ws.Range("K110").Value2 = 2010.6666
$ws.Range("M110").Value2 = 34.33339999999998
$ws.Range("H116").Value2 = 2082.2222
$ws.Range("I116").Value2 = 1330.5834
$ws.Range("K116").Value2 = 1330.5834
$ws.Range("M116").Value2 = 963.4166
$ws.Range("H122").Value2 = 12
$ws.Range("I122").Value2 = 12
$ws.Range("K122").Value2 = 36
$ws.Range("M122").Value2 = 2414
$ws.Range("H132").Value2 = 2500
$ws.Range("I132").Value2 = 2500
$ws.Range("K132").Value2 = 7500
$ws.Range("M132").Value2 = -4970

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 2082.2222
$ws.Range("I3").Value2 = 1330.5834
$ws.Range("K3").Value2 = 1330.5834
$ws.Range("M3").Value2 = -1216.5834
$ws.Range("H35").Value2 = 14998.75
$ws.Range("J35").Value2 = 14998.75
$ws.Range("L35").Value2 = 14998.75
$ws.Range("N35").Value2 = -15618.75
$ws.Range("H82").Value2 = 20034.916
$ws.Range("J82").Value2 = 40000
$ws.Range("L82").Value2 = 40000
$ws.Range("N82").Value2 = -40766
$ws.Range("H85").Value2 = 20034.916
$ws.Range("J85").Value2 = 40000
$ws.Range("L85").Value2 = 40000
$ws.Range("N85").Value2 = -42652

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 7064
$ws.Range("J16").Value2 = 10350
$ws.Range("L16").Value2 = 10350
$ws.Range("N16").Value2 = -10924
$ws.Range("H41").Value2 = 15925
$ws.Range("J41").Value2 = 16890
$ws.Range("L41").Value2 = 16890
$ws.Range("N41").Value2 = -17746
$ws.Range("H58").Value2 = 6203.579
$ws.Range("I58").Value2 = 5302.3125
$ws.Range("J58").Value2 = 11010.333
$ws.Range("K58").Value2 = 5302.3125
$ws.Range("L58").Value2 = 11010.333
$ws.Range("M58").Value2 = -5099.3125
$ws.Range("N58").Value2 = -11416.333
$ws.Range("H59").Value2 = 29977.375
$ws.Range("J59").Value2 = 34997.777
$ws.Range("L59").Value2 = 34997.777
$ws.Range("N59").Value2 = -37287.777
$ws.Range("H99").Value2 = 6599.6665
$ws.Range("I99").Value2 = 8900
$ws.Range("K99").Value2 = 8900
$ws.Range("M99").Value2 = -7402
$ws.Range("H113").Value2 = 7064
$ws.Range("J113").Value2 = 10350
$ws.Range("L113").Value2 = 10350
$ws.Range("N113").Value2 = -14690
$ws.Range("H126").Value2 = 6599.6665
$ws.Range("I126").Value2 = 8900
$ws.Range("K126").Value2 = 26700
$ws.Range("M126").Value2 = -24230
$ws.Range("H136").Value2 = 6203.579
$ws.Range("I136").Value2 = 5302.3125
$ws.Range("J136").Value2 = 11010.333
$ws.Range("K136").Value2 = 15906.9375
$ws.Range("L136").Value2 = 33030.999
$ws.Range("M136").Value2 = -13356.9375
$ws.Range("N136").Value2 = -38130.999
$ws.Range("H138").Value2 = 49999
$ws.Range("J138").Value2 = 49999
$ws.Range("L138").Value2 = 49999
$ws.Range("N138").Value2 = -60279
$ws.Range("H140").Value2 = 99994.5
$ws.Range("J140").Value2 = 99994.5
$ws.Range("L140").Value2 = 99994.5
$ws.Range("N140").Value2 = -110354.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value2 = 200
$ws.Range("I114").Value2 = 200
$ws.Range("K114").Value2 = 600
$ws.Range("M114").Value2 = 2654
$ws.Range("H139").Value2 = 1240.9166
$ws.Range("I139").Value2 = 899.1818
$ws.Range("K139").Value2 = 2697.5454
$ws.Range("M139").Value2 = 2442.4546

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value2 = 24095.875
$ws.Range("I43").Value2 = 1000
$ws.Range("J43").Value2 = 27395.285
$ws.Range("K43").Value2 = 1000
$ws.Range("L43").Value2 = 27395.285
$ws.Range("M43").Value2 = -849
$ws.Range("N43").Value2 = -27697.285
$ws.Range("H102").Value2 = 1577.2727
$ws.Range("I102").Value2 = 1744.4445
$ws.Range("K102").Value2 = 1744.4445
$ws.Range("M102").Value2 = -122.4445000000001
$ws.Range("H107").Value2 = 280.16666
$ws.Range("I107").Value2 = 263.33334
$ws.Range("K107").Value2 = 263.33334
$ws.Range("M107").Value2 = 1656.66666

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value2 = 4431.4
$ws.Range("I40").Value2 = 4663.154
$ws.Range("K40").Value2 = 4663.154
$ws.Range("M40").Value2 = -4527.154
$ws.Range("H61").Value2 = 3843.5
$ws.Range("I61").Value2 = 3843.5
$ws.Range("J61").Value2 = 0
$ws.Range("K61").Value2 = 3843.5
$ws.Range("L61").Value2 = 0
$ws.Range("M61").Value2 = -3641.5
$ws.Range("N61").ClearContents()
$ws.Range("H76").Value2 = 2250.75
$ws.Range("J76").Value2 = 2250.75
$ws.Range("L76").Value2 = 2250.75
$ws.Range("N76").Value2 = -2926.75
$ws.Range("H79").Value2 = 2250.75
$ws.Range("J79").Value2 = 2250.75
$ws.Range("L79").Value2 = 2250.75
$ws.Range("N79").Value2 = -4590.75
$ws.Range("H113").Value2 = 3843.5
$ws.Range("I113").Value2 = 3843.5
$ws.Range("J113").Value2 = 0
$ws.Range("K113").Value2 = 3843.5
$ws.Range("L113").Value2 = 0
$ws.Range("M113").Value2 = -1673.5
$ws.Range("N113").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value2 = 324.54544
$ws.Range("I107").Value2 = 310
$ws.Range("K107").Value2 = 930
$ws.Range("M107").Value2 = 990
$ws.Range("H110").Value2 = 55701.668
$ws.Range("J110").Value2 = 55701.668
$ws.Range("L110").Value2 = 55701.668
$ws.Range("N110").Value2 = -63881.668
$ws.Range("H126").Value2 = 5413.7144
$ws.Range("I126").Value2 = 5974.5
$ws.Range("K126").Value2 = 17923.5
$ws.Range("M126").Value2 = -15453.5
$ws.Range("H132").Value2 = 2781.7144
$ws.Range("I132").Value2 = 2781.7144
$ws.Range("J132").Value2 = 0
$ws.Range("K132").Value2 = 8345.143199999999
$ws.Range("L132").Value2 = 0
$ws.Range("M132").Value2 = -5815.143199999999
$ws.Range("N132").ClearContents()
